$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D3").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D4").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D5").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D6").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D7").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D8").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D9").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D10").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D11").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D12").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D13").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D14").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D15").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D16").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D17").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D18").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D19").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D20").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D21").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D22").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D23").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D24").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D25").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D26").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D27").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D28").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D29").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D30").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D31").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D32").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D33").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D34").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D35").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D36").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D37").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D38").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D39").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D40").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D41").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D42").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D43").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D44").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D45").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D46").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D47").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D48").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D49").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D50").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D51").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D52").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D53").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D54").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D55").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D56").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D57").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D58").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D59").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D60").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D61").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D62").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D63").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D64").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D65").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D66").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D67").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D68").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D69").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D70").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D71").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D72").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D73").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D74").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D75").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D76").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D77").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D78").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D79").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D80").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D81").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D82").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D83").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D84").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D85").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D86").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D87").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D88").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D89").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D90").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D91").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D92").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D93").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D94").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D95").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D96").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D97").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D98").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D99").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D100").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D101").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D102").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D103").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D104").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D105").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D106").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D107").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D108").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D109").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D110").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D111").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D112").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D113").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D114").Value = "2024-07-04T09:45:00.000Z"
$ws.Range("D115").Value = "2024-07-04T09:45:00.000Z"
